$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp shown in row 1
$ws.Range("A1").Value = "Datos actualizados a 11 de Agosto de 2020 a las 21:56"

# Helper-free, explicit per-row updates (country name in col A, stats in B:H)
# Each row is (RowNumber, Country, Total, NewCases, Active, Recovered, Critical, DeathsToday, Deaths)

# --- Row 4: Estados Unidos (name unchanged, stats refreshed) ---
$ws.Range("B4").Value = 5278963
$ws.Range("C4").Value = 27517
$ws.Range("D4").Value = 2722455
$ws.Range("E4").Value = 2389440
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 876
$ws.Range("H4").Value = 167068

# --- Row 22: Alemania (name unchanged, stats refreshed) ---
$ws.Range("B22").Value = 219409
$ws.Range("C22").Value = 909
$ws.Range("D22").Value = 198900
$ws.Range("E22").Value = 11242
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 9267

# --- Row 68: now Costa Rica (moved up from row 70, stats refreshed) ---
$ws.Range("A68").Value = "Costa Rica"
$ws.Range("B68").Value = 24508
$ws.Range("C68").Value = 636
$ws.Range("D68").Value = 7971
$ws.Range("E68").Value = 16282
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 11
$ws.Range("H68").Value = 255

# --- Row 69: now Etiopia (moved down from row 68, stats unchanged) ---
$ws.Range("A69").Value = "Etiopia"
$ws.Range("B69").Value = 24175
$ws.Range("C69").Value = 584
$ws.Range("D69").Value = 10696
$ws.Range("E69").Value = 13039
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 20
$ws.Range("H69").Value = 440

# --- Row 70: now Nepal (moved down from row 69, stats unchanged) ---
$ws.Range("A70").Value = "Nepal"
$ws.Range("B70").Value = 23948
$ws.Range("C70").Value = 638
$ws.Range("D70").Value = 16664
$ws.Range("E70").Value = 7201
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 4
$ws.Range("H70").Value = 83

# --- Row 137: Tunez (name unchanged, stats refreshed) ---
$ws.Range("B137").Value = 1738
$ws.Range("C137").Value = 21
$ws.Range("D137").Value = 1272
$ws.Range("E137").Value = 414
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 1
$ws.Range("H137").Value = 52

# --- Row 146: now Republica de Chipre (moved up from row 147, stats refreshed) ---
$ws.Range("A146").Value = "Republica de Chipre"
$ws.Range("B146").Value = 1277
$ws.Range("C146").Value = 25
$ws.Range("D146").Value = 870
$ws.Range("E146").Value = 387
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 1
$ws.Range("H146").Value = 20

# --- Row 147: now Georgia (moved down from row 146, stats unchanged) ---
$ws.Range("A147").Value = "Georgia"
$ws.Range("B147").Value = 1264
$ws.Range("C147").Value = 14
$ws.Range("D147").Value = 1054
$ws.Range("E147").Value = 193
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 17

# --- Row 161: now Aruba (moved up from row 164, stats refreshed) ---
$ws.Range("A161").Value = "Aruba"
$ws.Range("B161").Value = 717
$ws.Range("C161").Value = 87
$ws.Range("D161").Value = 114
$ws.Range("E161").Value = 600
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 3

# --- Row 162: now Crucero (moved down from row 161, stats unchanged) ---
$ws.Range("A162").Value = "Crucero"
$ws.Range("B162").Value = 712
$ws.Range("C162").Value = 0
$ws.Range("D162").Value = 651
$ws.Range("E162").Value = 48
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = 13

# --- Row 163: now San Marino (moved down from row 162, stats unchanged) ---
$ws.Range("A163").Value = "San Marino"
$ws.Range("B163").Value = 699
$ws.Range("C163").Value = 0
$ws.Range("D163").Value = 657
$ws.Range("E163").Value = 0
$ws.Range("F163").Value = 0
$ws.Range("G163").Value = 0
$ws.Range("H163").Value = 42

# --- Row 164: now Reunion (moved down from row 163, stats unchanged) ---
$ws.Range("A164").Value = "Reunion"
$ws.Range("B164").Value = 690
$ws.Range("C164").Value = 0
$ws.Range("D164").Value = 631
$ws.Range("E164").Value = 54
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 5

# --- Row 213: now Islas Malvinas (swapped with row 214) ---
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 13
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

# --- Row 214: now Montserrat (swapped with row 213) ---
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1
